{"js": "// Insert a new \"Body Text\" paragraph right after the figure-caption\n// paragraph (\"Verzi\u00f3k\u00f6vet\u00e9s alapok. ...\") and before the closing\n// paragraph (\"Ezek voltak az alapok, ...\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the caption paragraph (contains the picture + italic caption text).\nconst items = paragraphs.items;\nlet anchor = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Verzi\u00f3k\u00f6vet\u00e9s alapok.\") !== -1) {\n    anchor = items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph (image caption) not found.\");\n}\n\nconst newText =\n  \"A Git m\u0171k\u00f6d\u00e9s\u00e9t \u00fagyis elk\u00e9pzelhetj\u00fck, mintha Dropboxot/OneDrive-ot/Google drive-ot \" +\n  \"haszn\u00e1ln\u00e1nk annyi k\u00fcl\u00f6nbs\u00e9ggel, hogy itt a felh\u0151be val\u00f3 felt\u00f6lt\u00e9st mi ind\u00edtjuk k\u00e9zzel \" +\n  \"\u00e9s fejlettebb eszk\u00f6zeink vannak a verzi\u00f3k kezel\u00e9s\u00e9re \u00e9s a k\u00f6z\u00f6s munk\u00e1ra.\";\n\nconst newParagraph = anchor.insertParagraph(newText, Word.InsertLocation.after);\nnewParagraph.style = \"Body Text\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"Body Text\" paragraph right after the figure-caption\n# paragraph (\"Verzi\u00f3k\u00f6vet\u00e9s alapok. ...\") and before the closing\n# paragraph (\"Ezek voltak az alapok, ...\").\n\n$d = $word.ActiveDocument\n\n$newText = \"A Git m\u0171k\u00f6d\u00e9s\u00e9t \u00fagyis elk\u00e9pzelhetj\u00fck, mintha Dropboxot/OneDrive-ot/Google drive-ot haszn\u00e1ln\u00e1nk annyi k\u00fcl\u00f6nbs\u00e9ggel, hogy itt a felh\u0151be val\u00f3 felt\u00f6lt\u00e9st mi ind\u00edtjuk k\u00e9zzel \u00e9s fejlettebb eszk\u00f6zeink vannak a verzi\u00f3k kezel\u00e9s\u00e9re \u00e9s a k\u00f6z\u00f6s munk\u00e1ra.\"\n\n# Locate the paragraph that should immediately follow the new one.\n$rng = $d.Content\n$rng.Find.Execute(\"Ezek voltak az alapok\") | Out-Null\n\nif (-not $rng.Find.Found) {\n    throw \"Could not find the anchor paragraph ('Ezek voltak az alapok...').\"\n}\n\n# Expand the found range to the whole paragraph so InsertBefore lands\n# right at its start.\n$rng.Expand(4) | Out-Null\n\n# Insert the new paragraph's text followed by a paragraph mark; this\n# splits a brand-new paragraph in right before the found one, so it\n# inherits that (plain) paragraph formatting rather than the italic\n# caption formatting from the paragraph above.\n$rng.InsertBefore($newText + [char]13)\n\n# Locate the freshly split-off paragraph: it is the one immediately\n# before the \"Ezek voltak az alapok...\" paragraph.\n$followerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"Ezek voltak az alapok*\") {\n        $followerIndex = $i\n        break\n    }\n}\n\nif ($followerIndex -lt 2) {\n    throw \"Could not locate the newly inserted paragraph.\"\n}\n\n$newPara = $d.Paragraphs.Item($followerIndex - 1)\n$newPara.Style = \"Body Text\"\n\nWrite-Output \"Inserted paragraph. Document now has $($d.Paragraphs.Count) paragraphs.\"\n"}
